$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 403.33334
$ws.Range("I11").Value = 403.33334
$ws.Range("K11").Value = 403.33334
$ws.Range("M11").Value = -263.33334
$ws.Range("H18").Value = 2587.5
$ws.Range("I18").Value = 2925
$ws.Range("J18").Value = 2250
$ws.Range("K18").Value = 2925
$ws.Range("L18").Value = 2250
$ws.Range("M18").Value = -2641
$ws.Range("N18").Value = -2818
$ws.Range("H29").Value = 2317.1667
$ws.Range("I29").Value = 1003
$ws.Range("J29").Value = 2580
$ws.Range("K29").Value = 3009
$ws.Range("L29").Value = 7740
$ws.Range("M29").Value = -2728
$ws.Range("N29").Value = -8302
$ws.Range("H33").Value = 432.83334
$ws.Range("I33").Value = 447.30768
$ws.Range("J33").Value = 395.2
$ws.Range("K33").Value = 447.30768
$ws.Range("L33").Value = 395.2
$ws.Range("M33").Value = -218.30768
$ws.Range("N33").Value = -853.2
$ws.Range("H40").Value = 1852.7368
$ws.Range("I40").Value = 3262.5
$ws.Range("J40").Value = 1476.8
$ws.Range("K40").Value = 3262.5
$ws.Range("L40").Value = 1476.8
$ws.Range("M40").Value = -3087.5
$ws.Range("N40").Value = -1826.8
$ws.Range("H41").Value = 23811216
$ws.Range("I41").Value = 47620650
$ws.Range("J41").Value = 1784.2858
$ws.Range("K41").Value = 47620650
$ws.Range("L41").Value = 1784.2858
$ws.Range("M41").Value = -47620210
$ws.Range("N41").Value = -2664.2858
$ws.Range("H74").Value = 3933.3333
$ws.Range("I74").Value = 3900
$ws.Range("K74").Value = 3900
$ws.Range("M74").Value = -2964
$ws.Range("H77").Value = 3933.3333
$ws.Range("I77").Value = 3900
$ws.Range("K77").Value = 19500
$ws.Range("M77").Value = -14820
$ws.Range("H88").Value = 1237147.4
$ws.Range("I88").Value = 897.5
$ws.Range("J88").Value = 1546209.9
$ws.Range("K88").Value = 897.5
$ws.Range("L88").Value = 1546209.9
$ws.Range("M88").Value = -491.5
$ws.Range("N88").Value = -1547021.9
$ws.Range("H91").Value = 1237147.4
$ws.Range("I91").Value = 897.5
$ws.Range("J91").Value = 1546209.9
$ws.Range("K91").Value = 897.5
$ws.Range("L91").Value = 1546209.9
$ws.Range("M91").Value = 506.5
$ws.Range("N91").Value = -1549017.9
$ws.Range("H113").Value = 2659.9
$ws.Range("I113").Value = 2050.75
$ws.Range("J113").Value = 3066
$ws.Range("K113").Value = 2050.75
$ws.Range("L113").Value = 3066
$ws.Range("M113").Value = 1203.25
$ws.Range("N113").Value = -9574

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1536.0769
$ws.Range("I61").Value = 1274.3334
$ws.Range("K61").Value = 1274.3334
$ws.Range("M61").Value = -1062.3334
$ws.Range("H102").Value = 16677802
$ws.Range("I102").Value = 16677802
$ws.Range("K102").Value = 16677802
$ws.Range("M102").Value = -16676180
$ws.Range("H136").Value = 1536.0769
$ws.Range("I136").Value = 1274.3334
$ws.Range("K136").Value = 3823.0002
$ws.Range("M136").Value = -1273.0002

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2536.8572
$ws.Range("I20").Value = 2456.3333
$ws.Range("J20").Value = 2681.8
$ws.Range("K20").Value = 2456.3333
$ws.Range("L20").Value = 2681.8
$ws.Range("M20").Value = -2209.3333
$ws.Range("N20").Value = -3175.8
$ws.Range("H27").Value = 20542
$ws.Range("J27").Value = 20542
$ws.Range("L27").Value = 20542
$ws.Range("N27").Value = -20926
$ws.Range("H29").Value = 800
$ws.Range("I29").Value = 800
$ws.Range("K29").Value = 800
$ws.Range("M29").Value = -511
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H54").Value = 2317
$ws.Range("I54").Value = 731.625
$ws.Range("K54").Value = 731.625
$ws.Range("M54").Value = -247.625
$ws.Range("H107").Value = 1679.091
$ws.Range("I107").Value = 1175.1538
$ws.Range("J107").Value = 2407
$ws.Range("K107").Value = 1175.1538
$ws.Range("L107").Value = 2407
$ws.Range("M107").Value = 744.8462
$ws.Range("N107").Value = -6247

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 90910056
$ws.Range("I16").Value = 111112010
$ws.Range("K16").Value = 111112010
$ws.Range("M16").Value = -111111723
$ws.Range("H22").Value = 487.5
$ws.Range("I22").Value = 487.5
$ws.Range("K22").Value = 487.5
$ws.Range("M22").Value = -137.5
$ws.Range("H62").Value = 5717269.5
$ws.Range("I62").Value = 3013.2727
$ws.Range("K62").Value = 3013.2727
$ws.Range("M62").Value = -2389.2727
$ws.Range("H65").Value = 5717269.5
$ws.Range("I65").Value = 3013.2727
$ws.Range("K65").Value = 15066.3635
$ws.Range("M65").Value = -11946.3635
$ws.Range("H74").Value = 28925.4
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 28925.4
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 28925.4
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -30673.4
$ws.Range("H77").Value = 28925.4
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 28925.4
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 86776.20000000001
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -95512.20000000001
$ws.Range("H94").Value = 1349.4615
$ws.Range("I94").Value = 1221
$ws.Range("K94").Value = 1221
$ws.Range("M94").Value = -770
$ws.Range("H113").Value = 90910056
$ws.Range("I113").Value = 111112010
$ws.Range("K113").Value = 111112010
$ws.Range("M113").Value = -111109840

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2159.5881
$ws.Range("I69").Value = 999
$ws.Range("J69").Value = 2232.125
$ws.Range("K69").Value = 2997
$ws.Range("L69").Value = 6696.375
$ws.Range("M69").Value = -2186
$ws.Range("N69").Value = -8318.375
$ws.Range("H72").Value = 2159.5881
$ws.Range("I72").Value = 999
$ws.Range("J72").Value = 2232.125
$ws.Range("K72").Value = 8991
$ws.Range("L72").Value = 20089.125
$ws.Range("M72").Value = -4935
$ws.Range("N72").Value = -28201.125
$ws.Range("H137").Value = 17081.223
$ws.Range("I137").Value = 5500
$ws.Range("J137").Value = 20390.143
$ws.Range("K137").Value = 16500
$ws.Range("L137").Value = 61170.429
$ws.Range("M137").Value = -11400
$ws.Range("N137").Value = -71370.429

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18754942
$ws.Range("I70").Value = 17861452
$ws.Range("J70").Value = 20005830
$ws.Range("K70").Value = 17861452
$ws.Range("L70").Value = 20005830
$ws.Range("M70").Value = -17861182
$ws.Range("N70").Value = -20006370
$ws.Range("H73").Value = 18754942
$ws.Range("I73").Value = 17861452
$ws.Range("J73").Value = 20005830
$ws.Range("K73").Value = 17861452
$ws.Range("L73").Value = 20005830
$ws.Range("M73").Value = -17860516
$ws.Range("N73").Value = -20007702
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -46134
$ws.Range("H118").Value = 17333.334
$ws.Range("J118").Value = 17333.334
$ws.Range("L118").Value = 17333.334
$ws.Range("N118").Value = -20647.334
$ws.Range("H122").Value = 1350
$ws.Range("I122").Value = 1350
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4050
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1600
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3457.5334
$ws.Range("I126").Value = 1889.3334
$ws.Range("J126").Value = 4503
$ws.Range("K126").Value = 5668.0002
$ws.Range("L126").Value = 13509
$ws.Range("M126").Value = -3198.0002
$ws.Range("N126").Value = -18449
$ws.Range("H132").Value = 2387.1924
$ws.Range("I132").Value = 2178.3684
$ws.Range("J132").Value = 2954
$ws.Range("K132").Value = 6535.1052
$ws.Range("L132").Value = 8862
$ws.Range("M132").Value = -4005.1052
$ws.Range("N132").Value = -13922

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2958.2
$ws.Range("I40").Value = 2677.2
$ws.Range("K40").Value = 2677.2
$ws.Range("M40").Value = -2541.2
$ws.Range("H46").Value = 4816.7144
$ws.Range("I46").Value = 1250.2
$ws.Range("J46").Value = 5931.25
$ws.Range("K46").Value = 1250.2
$ws.Range("L46").Value = 5931.25
$ws.Range("M46").Value = -1062.2
$ws.Range("N46").Value = -6307.25
$ws.Range("H122").Value = 27781264
$ws.Range("J122").Value = 4166.3335
$ws.Range("L122").Value = 12499.0005
$ws.Range("N122").Value = -17399.0005
$ws.Range("H136").Value = 10219.546
$ws.Range("I136").Value = 10219.546
$ws.Range("K136").Value = 30658.638
$ws.Range("M136").Value = -28108.638

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15299909
$ws.Range("I122").Value = 17339496
$ws.Range("J122").Value = 3002.5
$ws.Range("K122").Value = 52018488
$ws.Range("L122").Value = 9007.5
$ws.Range("M122").Value = -52016038
$ws.Range("N122").Value = -13907.5
$ws.Range("H126").Value = 76924530
$ws.Range("I126").Value = 76924530
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 230773590
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -230771120
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 738.93335
$ws.Range("I136").Value = 381.91666
$ws.Range("J136").Value = 2167
$ws.Range("K136").Value = 1145.74998
$ws.Range("L136").Value = 6501
$ws.Range("M136").Value = 1404.25002
$ws.Range("N136").Value = -11601
